# Adds Word's automatic "possible spelling error" markers (<w:proofErr/>)
# around the Latin-script / non-dictionary words Otchet, analiz and гитхаб,
# as happens when Word's background spell checker re-scans the document
# (e.g. after the text was retyped/pasted), splitting the runs that carry
# those words so the proofErr start/end bracket them exactly.

$d = $word.ActiveDocument

# Namespace-qualified "Flat OPC" wrapper required by this host's
# Range.InsertXML: it parses the payload as a full single-part package and
# splices the <w:body> content of the wrapping <w:p> into the target Range,
# so it must be used with a Range that excludes the paragraph mark (the
# Range is otherwise swallowed whole, paragraph mark and all).
function New-FlatOpcParagraph([string]$innerXml) {
    return '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $innerXml + '</w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParagraphRuns([string]$containsText, [string]$innerXml) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -match [regex]::Escape($containsText)) {
            $start = $p.Range.Start
            $end = $p.Range.End - 1   # exclude the paragraph mark
            $r = $d.Range($start, $end)
            $r.InsertXML((New-FlatOpcParagraph $innerXml))
            return $true
        }
    }
    return $false
}

# 1) "Создание файла Otchet и добавление коммитов к analiz."
#    -> wrap "Otchet" and "analiz" with spellStart/spellEnd proofErr marks.
$inner1 = (
    '<w:r><w:t xml:space="preserve">Создание файла </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Otchet</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r w:rsidRPr="001D6867"><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">и добавление коммитов к </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>analiz</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r w:rsidRPr="001D6867"><w:t>.</w:t></w:r>'
)
Set-ParagraphRuns "Создание файла" $inner1 | Out-Null

# 2) "Изменения в файл Otchet. " -> split the "Otchet. " run so
#    spellStart/spellEnd bracket just "Otchet".
$inner2 = (
    '<w:r><w:t xml:space="preserve">Изменения в файл </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Otchet</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r>'
)
Set-ParagraphRuns "Изменения в файл" $inner2 | Out-Null

# 3) "Загрузка на гитхаб" -> split off "гитхаб" and bracket it.
$inner3 = (
    '<w:r><w:t xml:space="preserve">Загрузка на </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>гитхаб</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)
Set-ParagraphRuns "Загрузка на" $inner3 | Out-Null

Write-Output "done"
